# chore: Fix blueprint parameter sorting
# Adds a new "DayNightCycle" worksheet (after "Animals") containing the
# day-cycle parameter table (IDS / Duration / Name columns).

$wb = $excel.ActiveWorkbook
$animals = $wb.Worksheets.Item(1)

# Insert the new sheet right after "Animals" so sheet order/ids match
# (Animals=sheetId 1, DayNightCycle=sheetId 2).
$ws = $wb.Worksheets.Add($null, $animals)
$ws.Name = "DayNightCycle"

# Header row
$ws.Range("A1").Value = "IDS"
$ws.Range("B1").Value = "Duration"
$ws.Range("C1").Value = "Name"

# Data rows: Id, Duration (hours), Name
$ws.Range("A2").Value = "Mañana"
$ws.Range("B2").Value = 24
$ws.Range("C2").Value = "Manana"

$ws.Range("A3").Value = "Mediodía"
$ws.Range("B3").Value = 24
$ws.Range("C3").Value = "Mediodia"

$ws.Range("A4").Value = "Tarde"
$ws.Range("B4").Value = 24
$ws.Range("C4").Value = "Tarde"

$ws.Range("A5").Value = "Atardecer"
$ws.Range("B5").Value = 24
$ws.Range("C5").Value = "Atardecer"

$ws.Range("A6").Value = "Anochecer"
$ws.Range("B6").Value = 24
$ws.Range("C6").Value = "Anochecer"

$ws.Range("A7").Value = "Madrugada"
$ws.Range("B7").Value = 24
$ws.Range("C7").Value = "Madrugada"

# Match the look & feel of the existing "Animals" sheet (same font/style
# family) by copying its cell formats onto the new table.
$animals.Range("A1").Copy()
$ws.Range("A1:C7").PasteSpecial(-4122) # xlPasteFormats

# Duration (numbers) and Name columns are right aligned in the source data.
$ws.Range("B2:B7").HorizontalAlignment = -4152 # xlRight
$ws.Range("C2:C7").HorizontalAlignment = -4152 # xlRight

# Keep "Animals" as the active sheet, same as before the edit.
$animals.Activate()
